$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Existing data runs from row 2 (first day) through row 126 (most recent
# day, serial date 45682). Append two more daily rows (127, 128) for the
# next two days, carrying forward the same metric values (columns B..J)
# that row 126 already has, and the same date-cell formatting used by the
# rest of column A.

$sourceRow = 126
$newRows = @(127, 128)

foreach ($newRow in $newRows) {
    # Column A: copy the date cell's formatting from the row above, then
    # set the new (incremented) date serial value.
    $ws.Cells.Item($sourceRow, 1).Copy()
    $ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
    $ws.Cells.Item($newRow, 1).Value = $ws.Cells.Item($sourceRow, 1).Value2 + ($newRow - $sourceRow)

    # Columns B..J: plain numeric values, copied from row 126.
    for ($col = 2; $col -le 10; $col++) {
        $ws.Cells.Item($newRow, $col).Value = $ws.Cells.Item($sourceRow, $col).Value2
    }
}

$excel.CutCopyMode = $false
